
$wb = $excel.ActiveWorkbook

# --- Rename the original "Acc_Upfront" sheet to "Acc_Upfront1" ---
$ws1 = $wb.Worksheets.Item("Acc_Upfront")
$ws1.Name = "Acc_Upfront1"

# --- Insert two brand-new sheets right after it: Acc_Upfront2, Acc_Upfront3 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Acc_Upfront2"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Acc_Upfront3"

# --- Populate Acc_Upfront2 with the "2949 / 2950" record group (old rows 5:6) ---
$ws1.Range("A1:I1").Copy($ws2.Range("A1:I1"))
$ws1.Range("A5:I6").Copy($ws2.Range("A2:I3"))

# --- Populate Acc_Upfront3 with the "3820 / 3821" record group (old rows 8:9) ---
$ws1.Range("A1:I1").Copy($ws3.Range("A1:I1"))
$ws1.Range("A8:I9").Copy($ws3.Range("A2:I3"))

# --- Trim Acc_Upfront1 back down to just its own group (header + rows 2:3) ---
$ws1.Rows("4:9").Delete()

# --- Fix up each sheet's selection / active cell ---
$null = $ws1.Range("E2").Select()
$null = $ws2.Range("E3").Select()
$null = $ws3.Range("F3").Select()

# --- Acc_Upfront3 is the active/selected tab now (was Acc_Upfront1 before) ---
$ws3.Activate()
